$wb = $excel.ActiveWorkbook

$wsAll  = $wb.Worksheets.Item("ALL")
$wsPlay = $wb.Worksheets.Item("Playable (untested)")

# --- Read the source row (the "airbustr" entry) from ALL!1528 ---
$srcRow = 1528
$lastCol = 13  # A..M
$vals = @{}
for ($c = 1; $c -le $lastCol; $c++) {
    $vals[$c] = $wsAll.Cells.Item($srcRow, $c).Value()
}

# --- Append it as a new row at the bottom of "Playable (untested)" ---
$destRow = 357
for ($c = 1; $c -le $lastCol; $c++) {
    $v = $vals[$c]
    if ($c -eq 1) {
        $wsPlay.Cells.Item($destRow, $c).Value2 = $destRow
    } elseif ($v -ne $null -and $v -ne "") {
        $wsPlay.Cells.Item($destRow, $c).Value2 = $v
    }
}
# New column N on the destination row
$wsPlay.Cells.Item($destRow, 14).Value2 = "tilemaps"

# --- Remove the now-duplicated row from ALL (shifts everything below up) ---
$wsAll.Rows.Item($srcRow).Delete()

# --- Fix up the AutoFilter range on ALL now that it has one fewer row ---
$wsAll.AutoFilterMode = $false
$wsAll.Range("A1:M1685").AutoFilter()

# --- Fix up the workbook-level defined names that pointed at ALL's old extent ---
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "ALL!_FilterDatabase") {
        $n.RefersTo = "=ALL!`$A`$1:`$M`$1685"
    } elseif ($n.Name -eq "ALL!LIST") {
        $n.RefersTo = "=ALL!`$B`$1:`$M`$1685"
    }
}

# --- Selection / active sheet bookkeeping ---
$wsPlay.Activate()
$wsPlay.Cells.Item(357, 14).Select()
